# Pandoc docx writer bugfix: avoid emitting duplicate <w:abstractNum> /
# <w:num> entries in word/numbering.xml. Reference.docx's sample list
# styles were being copied into the generated document alongside the
# autogenerated numbering entries, producing two <w:abstractNum
# w:abstractNumId="990"> elements and two <w:num w:numId="1000">
# elements. Desktop Word tolerated the clash but Word Online choked on
# it, so drop the extra (first/earlier) copy of each, keeping a single
# abstractNum/num pair.

$d = $word.ActiveDocument

# Word's WordOpenXML property round-trips the *entire* package
# (all parts, including word/numbering.xml) as a single flat-OPC XML
# document -- grab it, surgically remove the duplicated numbering
# elements, and write it back.
$xml = $d.WordOpenXML

# 1) Remove the first (duplicate) <w:abstractNum w:abstractNumId="990">
#    block -- the one pandoc's reference.docx contributed (nsid
#    170cd2de). The second abstractNum (nsid 2c1ae401) -- the one
#    pandoc actually generated for this document's list -- is left
#    alone.
$dupAbstractNum = '<w:abstractNum w:abstractNumId="990"><w:nsid w:val="170cd2de" /><w:multiLevelType w:val="multilevel" /><w:lvl w:ilvl="0"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="0" /></w:tabs><w:ind w:left="480" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="1"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="720" /></w:tabs><w:ind w:left="1200" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="2"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="1440" /></w:tabs><w:ind w:left="1920" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="3"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="2160" /></w:tabs><w:ind w:left="2640" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="4"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="2880" /></w:tabs><w:ind w:left="3360" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="5"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="3600" /></w:tabs><w:ind w:left="4080" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="6"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="4320" /></w:tabs><w:ind w:left="4800" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="7"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="5040" /></w:tabs><w:ind w:left="5520" w:hanging="480" /></w:pPr></w:lvl><w:lvl w:ilvl="8"><w:numFmt w:val="bullet" /><w:lvlText w:val=" " /><w:lvlJc w:val="left" /><w:pPr><w:tabs><w:tab w:val="num" w:pos="5760" /></w:tabs><w:ind w:left="6240" w:hanging="480" /></w:pPr></w:lvl></w:abstractNum>'

$xml = $xml.Replace($dupAbstractNum, "")

# 2) Remove one of the two identical <w:num w:numId="1000"> wrapper
#    elements (both point at abstractNumId 990) -- only one is needed.
$numEntry = '<w:num w:numId="1000"><w:abstractNumId w:val="990" /></w:num>'
$firstIdx = $xml.IndexOf($numEntry)
$xml = $xml.Remove($firstIdx, $numEntry.Length)

$d.WordOpenXML = $xml
